$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 7 new rows at row 2, shifting existing data (rows 2-21) down to rows 9-28
$ws.Rows.Item(2).Resize(7).Insert()
$ws.Range("A2:C8").ClearFormats()

# Populate the 7 newly inserted rows (2-8) with new data
$newTop = New-Object "object[,]" 7,3
$newTop[0,0] = -0.0186313893646001
$newTop[0,1] = -0.1148426681756973
$newTop[0,2] = 0.057115901261568
$newTop[1,0] = 0.0178678091615438
$newTop[1,1] = -0.1070541366934776
$newTop[1,2] = 0.1860084682703018
$newTop[2,0] = -0.007941247895359899
$newTop[2,1] = 0.0074830991216003
$newTop[2,2] = 0.0245873257517814
$newTop[3,0] = -0.039248090237379
$newTop[3,1] = -0.008552113547921099
$newTop[3,2] = 0.0209221355617046
$newTop[4,0] = 0.00167987938039
$newTop[4,1] = -0.0282525178045034
$newTop[4,2] = 0.0545197241008281
$newTop[5,0] = 0.00167987938039
$newTop[5,1] = -0.0485637858510017
$newTop[5,2] = 0.0387899428606033
$newTop[6,0] = 0.0103847095742821
$newTop[6,1] = 0.0097738439217209
$newTop[6,2] = -0.0050396383740007
$ws.Range("A2:C8").Value = $newTop

# Append a new row 29 with new data
$newBottom = New-Object "object[,]" 1,3
$newBottom[0,0] = -0.0229074470698833
$newBottom[0,1] = 0.0600175112485885
$newBottom[0,2] = 0.0545197241008281
$ws.Range("A29:C29").Value = $newBottom

